# Append 9 new rows (regcntr_id/machine_id pairs 10002-10010 / 10021-10029) to
# the "master-reg_center_machine" master data sheet, matching the existing
# row pattern (lang_code="eng", is_active=TRUE, cr_by="superadmin",
# cr_dtimes="now()").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(10002, 10021),
    @(10003, 10022),
    @(10004, 10023),
    @(10005, 10024),
    @(10006, 10025),
    @(10007, 10026),
    @(10008, 10027),
    @(10009, 10028),
    @(10010, 10029)
)

$startRow = 22
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $regcntrId = $newRows[$i][0]
    $machineId = $newRows[$i][1]

    $ws.Cells.Item($r, 1).Value = $regcntrId
    $ws.Cells.Item($r, 2).Value = $machineId
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
}

# Reproduce the final selection state left behind in the saved file: the
# cursor rests on the row below the new data, with entire rows selected
# down to the bottom of the sheet.
$lastDataRow = $startRow + $newRows.Count - 1
$ws.Range("A" + ($lastDataRow + 1) + ":XFD1048576").Select() | Out-Null

# Page setup was touched (portrait orientation, 300 dpi) before saving.
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PrintQuality = @(300, 300)
